$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Bring the "git / line ending" and "Remote / My own git server" rows
# (currently rows 47-48) up into rows 37-38, carrying their values and
# cell formatting, overwriting the old "Markdown" cheat-sheet rows that
# used to live there.
$ws.Range("A47:C48").Copy($ws.Range("A37"))

# The two relocated rows get their own explicit heights.
$ws.Rows(37).RowHeight = 174.75
$ws.Rows(38).RowHeight = 120

# Remove the rest of the old Markdown section (rows 39-46) plus the
# now-duplicated source rows (47-48).
$ws.Rows("39:48").Delete()

# Re-create row 47 as a blank, formatted-only row (matches the leftover
# empty row left behind in the sheet).
$ws.Rows(47).RowHeight = 92.25

$ws.Range("A36").Select()
$ws.Application.ActiveWindow.ScrollRow = 36
$ws.Range("C38").Select()
